$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 10000570
$ws.Range("I33").Value = 16667098
$ws.Range("J33").Value = 777.8
$ws.Range("K33").Value = 16667098
$ws.Range("L33").Value = 777.8
$ws.Range("M33").Value = -16666869
$ws.Range("N33").Value = -1235.8
$ws.Range("H98").Value = 14337
$ws.Range("I98").Value = 3822.2273
$ws.Range("K98").Value = 3822.2273
$ws.Range("M98").Value = -2324.2273
$ws.Range("H112").Value = 37700.633
$ws.Range("I112").Value = 5207.143
$ws.Range("J112").Value = 113518.78
$ws.Range("K112").Value = 15621.429
$ws.Range("L112").Value = 340556.34
$ws.Range("M112").Value = -14513.429
$ws.Range("N112").Value = -342772.34
$ws.Range("H113").Value = 6519.6924
$ws.Range("I113").Value = 7790.5
$ws.Range("J113").Value = 4486.4
$ws.Range("K113").Value = 7790.5
$ws.Range("L113").Value = 4486.4
$ws.Range("M113").Value = -4536.5
$ws.Range("N113").Value = -10994.4
$ws.Range("H122").Value = 14337
$ws.Range("I122").Value = 3822.2273
$ws.Range("K122").Value = 11466.6819
$ws.Range("M122").Value = -9016.6819
$ws.Range("H132").Value = 2230188.8
$ws.Range("I132").Value = 1898.2778
$ws.Range("K132").Value = 5694.8334
$ws.Range("M132").Value = -3164.8334
$ws.Range("H135").Value = 4055.5715
$ws.Range("I135").Value = 4284.115
$ws.Range("K135").Value = 38557.035
$ws.Range("M135").Value = -36022.035
$ws.Range("H137").Value = 16199.105
$ws.Range("I137").Value = 5365.6665
$ws.Range("K137").Value = 16096.9995
$ws.Range("M137").Value = -13546.9995
$ws.Range("H138").Value = 3784.5898
$ws.Range("I138").Value = 4050.9167
$ws.Range("J138").Value = 3666.2222
$ws.Range("K138").Value = 12152.7501
$ws.Range("L138").Value = 10998.6666
$ws.Range("M138").Value = -7012.750100000001
$ws.Range("N138").Value = -21278.6666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1305.7407
$ws.Range("I45").Value = 1295.3
$ws.Range("K45").Value = 1295.3
$ws.Range("M45").Value = -918.3
$ws.Range("H61").Value = 823554.5
$ws.Range("I61").Value = 3928
$ws.Range("J61").Value = 3692247
$ws.Range("K61").Value = 3928
$ws.Range("L61").Value = 3692247
$ws.Range("M61").Value = -3716
$ws.Range("N61").Value = -3692671
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H94").Value = 29999.5
$ws.Range("J94").Value = 29999.5
$ws.Range("L94").Value = 29999.5
$ws.Range("N94").Value = -31801.5
$ws.Range("H97").Value = 4277.375
$ws.Range("I97").Value = 4277.375
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 4277.375
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -3781.375
$ws.Range("N97").ClearContents()
$ws.Range("H124").Value = 25000
$ws.Range("J124").Value = 25000
$ws.Range("L124").Value = 25000
$ws.Range("N124").Value = -34820
$ws.Range("H125").Value = 75554.664
$ws.Range("J125").Value = 75554.664
$ws.Range("L125").Value = 75554.664
$ws.Range("N125").Value = -85394.664
$ws.Range("H132").Value = 833370.6
$ws.Range("I132").Value = 3718.3572
$ws.Range("K132").Value = 11155.0716
$ws.Range("M132").Value = -8625.071599999999
$ws.Range("H136").Value = 823554.5
$ws.Range("I136").Value = 3928
$ws.Range("J136").Value = 3692247
$ws.Range("K136").Value = 11784
$ws.Range("L136").Value = 11076741
$ws.Range("M136").Value = -9234
$ws.Range("N136").Value = -11081841

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 22342.805
$ws.Range("I20").Value = 6765.5454
$ws.Range("J20").Value = 40379.633
$ws.Range("K20").Value = 6765.5454
$ws.Range("L20").Value = 40379.633
$ws.Range("M20").Value = -6518.5454
$ws.Range("N20").Value = -40873.633
$ws.Range("H94").Value = 1872.8
$ws.Range("I94").Value = 1872.8
$ws.Range("K94").Value = 1872.8
$ws.Range("M94").Value = -1421.8
$ws.Range("H99").Value = 3576.1765
$ws.Range("I99").Value = 2710.5
$ws.Range("K99").Value = 2710.5
$ws.Range("M99").Value = -1212.5
$ws.Range("H105").Value = 2887.7693
$ws.Range("I105").Value = 2589.0952
$ws.Range("J105").Value = 4142.2
$ws.Range("K105").Value = 2589.0952
$ws.Range("L105").Value = 4142.2
$ws.Range("M105").Value = -842.0952000000002
$ws.Range("N105").Value = -7636.2
$ws.Range("H134").Value = 32611.875
$ws.Range("I134").Value = 35218.43
$ws.Range("J134").Value = 30584.555
$ws.Range("K134").Value = 105655.29
$ws.Range("L134").Value = 91753.66500000001
$ws.Range("M134").Value = -103120.29
$ws.Range("N134").Value = -96823.66500000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1468.0625
$ws.Range("I22").Value = 649.375
$ws.Range("J22").Value = 2286.75
$ws.Range("K22").Value = 649.375
$ws.Range("L22").Value = 2286.75
$ws.Range("M22").Value = -299.375
$ws.Range("N22").Value = -2986.75
$ws.Range("H33").Value = 999.5
$ws.Range("I33").Value = 999.5
$ws.Range("K33").Value = 999.5
$ws.Range("M33").Value = -620.5
$ws.Range("H36").Value = 6500
$ws.Range("I36").Value = 5000
$ws.Range("J36").Value = 8000
$ws.Range("K36").Value = 5000
$ws.Range("L36").Value = 8000
$ws.Range("M36").Value = -4612
$ws.Range("N36").Value = -8776
$ws.Range("H40").Value = 6500
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -4840
$ws.Range("N40").Value = -8320
$ws.Range("H134").Value = 6336.273
$ws.Range("I134").Value = 1875.4286
$ws.Range("J134").Value = 100014
$ws.Range("K134").Value = 5626.2858
$ws.Range("L134").Value = 300042
$ws.Range("M134").Value = -3091.2858
$ws.Range("N134").Value = -305112

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 5750.846
$ws.Range("J39").Value = 7110.6665
$ws.Range("L39").Value = 21331.9995
$ws.Range("N39").Value = -21919.9995
$ws.Range("H106").Value = 3600
$ws.Range("I106").Value = 1200
$ws.Range("K106").Value = 3600
$ws.Range("M106").Value = -2654
$ws.Range("H107").Value = 3228.65
$ws.Range("I107").Value = 607.9167
$ws.Range("J107").Value = 7159.75
$ws.Range("K107").Value = 1823.7501
$ws.Range("L107").Value = 21479.25
$ws.Range("M107").Value = 96.24990000000003
$ws.Range("N107").Value = -25319.25
$ws.Range("H113").Value = 966.6667
$ws.Range("J113").Value = 1100
$ws.Range("L113").Value = 3300
$ws.Range("N113").Value = -7640

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2199.6365
$ws.Range("I80").Value = 1919.6
$ws.Range("K80").Value = 1919.6
$ws.Range("M80").Value = -921.5999999999999
$ws.Range("H83").Value = 2199.6365
$ws.Range("I83").Value = 1919.6
$ws.Range("K83").Value = 9598
$ws.Range("M83").Value = -4606
$ws.Range("H107").Value = 502.8889
$ws.Range("I107").Value = 502.8889
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 502.8889
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1417.1111
$ws.Range("N107").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4229.222
$ws.Range("I40").Value = 3353.4
$ws.Range("K40").Value = 3353.4
$ws.Range("M40").Value = -3217.4
$ws.Range("H45").Value = 9000
$ws.Range("I45").Value = 9000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 9000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -8593
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 3967.7778
$ws.Range("I61").Value = 2894.7334
$ws.Range("K61").Value = 2894.7334
$ws.Range("M61").Value = -2692.7334
$ws.Range("H68").Value = 7990
$ws.Range("I68").Value = 6625
$ws.Range("J68").Value = 8900
$ws.Range("K68").Value = 6625
$ws.Range("L68").Value = 8900
$ws.Range("M68").Value = -5876
$ws.Range("N68").Value = -10398
$ws.Range("H71").Value = 7990
$ws.Range("I71").Value = 6625
$ws.Range("J71").Value = 8900
$ws.Range("K71").Value = 33125
$ws.Range("L71").Value = 44500
$ws.Range("M71").Value = -29381
$ws.Range("N71").Value = -51988
$ws.Range("H93").Value = 11406.692
$ws.Range("J93").Value = 2938.5
$ws.Range("N93").Value = -5434.5
$ws.Range("H113").Value = 3967.7778
$ws.Range("I113").Value = 2894.7334
$ws.Range("K113").Value = 2894.7334
$ws.Range("M113").Value = -724.7334000000001
$ws.Range("H136").Value = 1457508.2
$ws.Range("I136").Value = 24684.666
$ws.Range("J136").Value = 2449463
$ws.Range("K136").Value = 74053.99800000001
$ws.Range("L136").Value = 7348389
$ws.Range("M136").Value = -71503.99800000001
$ws.Range("N136").Value = -7353489

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 481.81818
$ws.Range("I107").Value = 469
$ws.Range("K107").Value = 1407
$ws.Range("M107").Value = 513
$ws.Range("H122").Value = 3193.879
$ws.Range("I122").Value = 2974.8462
$ws.Range("K122").Value = 8924.5386
$ws.Range("M122").Value = -6474.5386
$ws.Range("H132").Value = 305601.03
$ws.Range("I132").Value = 7778.8
$ws.Range("K132").Value = 23336.4
$ws.Range("M132").Value = -20806.4
